$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2   = "249,250,248"
    3   = "250,249,250"
    4   = "244,246,251"
    5   = "251,251,246"
    6   = "252,250,245"
    7   = "251,253,248"
    8   = "237,240,242"
    9   = "231,236,239"
    10  = "209,216,211"
    11  = "192,199,211"
    12  = "193,202,216"
    13  = "184,192,206"
    14  = "185,187,169"
    15  = "154,163,183"
    16  = "130,141,150"
    29  = "84,101,165"
    30  = "96,119,226"
    31  = "109,138,210"
    32  = "123,153,220"
    33  = "128,192,214"
    34  = "120,162,192"
    35  = "141,169,198"
    36  = "217,235,254"
    37  = "123,142,202"
    44  = "193,201,232"
    45  = "104,108,142"
    47  = "182,189,174"
    48  = "174,177,151"
    49  = "145,150,150"
    51  = "113,170,137"
    52  = "164,179,166"
    53  = "171,183,150"
    54  = "151,181,120"
    55  = "134,153,109"
    59  = "208,217,185"
    60  = "109,127,117"
    62  = "215,224,144"
    64  = "255,254,178"
    72  = "253,240,150"
    73  = "229,219,117"
    74  = "253,236,133"
    75  = "251,229,115"
    76  = "255,246,186"
    89  = "254,208,164"
    90  = "255,242,220"
    91  = "255,225,214"
    105 = "135,102,42"
    106 = "192,166,102"
}

foreach ($row in $values.Keys) {
    $cell = $ws.Cells.Item($row, 10)
    $cell.NumberFormat = "@"
    $cell.Value = $values[$row]
}
